$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Fix two existing GERACAO values (H32, H69)
# ---------------------------------------------------------------
$ws.Range("H32").Value = "JUVENTUDE RELEVANTE MOÇAS"
$ws.Range("H69").Value = "MULHERES IECG"

# ---------------------------------------------------------------
# 2. Add the two new respondents as rows 73 and 74
# ---------------------------------------------------------------
# Row 73 - Ricardo Hildebrand Camargo
$ws.Range("A73").Value = "rhildebrandcamargo@gmail.com"
$ws.Range("B73").Value = "RICARDO HILDEBRAND CAMARGO"
$ws.Range("C73").Value = 69809500106
$ws.Range("D73").Value = 28908
$ws.Range("D71").Copy()
$ws.Range("D73").PasteSpecial(-4122)
$ws.Range("E73").Value = "Campo Grande - MS"
$ws.Range("F73").Value = "Rua Ouro Negro, 186`nCasa"
$ws.Range("F73").WrapText = $true
$ws.Range("G73").Value = 67992625238
$ws.Range("G71").Copy()
$ws.Range("G73").PasteSpecial(-4122)
$ws.Range("H73").Value = "HOMENS IECG"
$ws.Range("I73").Value = "Souza"
$ws.Range("J73").Value = "Bispo EVALDO"
$ws.Range("K73").Value = "TERÇA FEIRA - 19H - IECG CENTRO"
$ws.Rows.Item(73).RowHeight = 15.75

# Row 74 - Dayana Rocha da Silva
$ws.Range("A74").Value = "dayanarochadasilva@yahoo.com.br"
$ws.Range("B74").Value = "Dayana rocha da Silva"
$ws.Range("C74").Value = 71064141153
$ws.Range("D74").Value = 29995
$ws.Range("D71").Copy()
$ws.Range("D74").PasteSpecial(-4122)
$ws.Range("E74").Value = "Campo Grande - MS"
$ws.Range("F74").Value = "Rua ouro negro 186 vila Marcos Roberto"
$ws.Range("G74").Value = "67 99238-4001"
$ws.Range("G71").Copy()
$ws.Range("G74").PasteSpecial(-4122)
$ws.Range("H74").Value = "MULHERES IECG"
$ws.Range("I74").Value = "Mirtes"
$ws.Range("J74").Value = "Prª LIU"
$ws.Range("K74").Value = "TERÇA FEIRA - 19H - IECG CENTRO"
$ws.Rows.Item(74).RowHeight = 15.75

# ---------------------------------------------------------------
# 3. Extend the autofilter / filter-database range to cover the
#    newly added rows (A1:K72)
# ---------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:K72").AutoFilter() | Out-Null

foreach ($n in $wb.Names) {
    $n.RefersTo = "='Respostas ao formulário 1'!`$A`$1:`$K`$72"
}

# ---------------------------------------------------------------
# 4. Row heights for the header/first data block (rows 1-21)
# ---------------------------------------------------------------
$ws.Rows.Item("1:21").RowHeight = 12.75

# ---------------------------------------------------------------
# 5. Update frozen pane scroll position & selection
# ---------------------------------------------------------------
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("H73:K74").Select() | Out-Null

# ---------------------------------------------------------------
# 6. Page setup (paper size / orientation)
# ---------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
